$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D2:E51 range to Text format so numeric-looking strings
# (e.g. "25.872.40", "5.440") are preserved exactly as typed, not
# coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.872.40"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.631.92"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "214.47"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "0.5116"
$ws.Range("E6").Value = "  +1.32%  "
$ws.Range("D7").Value = "1.003"
$ws.Range("D8").Value = "0.2548"
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("D9").Value = "0.06337"
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").Value = "19.41"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").Value = "0.07748"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "4.262"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "1.639.85"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "0.5405"
$ws.Range("D15").Value = "0.0₅7702"
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("D16").Value = "63.91"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "25.881.39"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "4.417"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").Value = "194.47"
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("D21").Value = "9.899"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "6.009"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "1.005"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "1.857"
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("D25").Value = "140.77"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").Value = "0.1186"
$ws.Range("E26").Value = "  +4.57%  "
$ws.Range("D27").Value = "6.808"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "15.55"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").Value = "0.04897"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").Value = "3.239"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").Value = "3.154"
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("D33").Value = "1.525"
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("D34").Value = "2.368"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").Value = "0.8870"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").Value = "2.572"
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("D37").Value = "1.138.35"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("D38").Value = "0.5381"
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("D39").Value = "0.01544"
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("D40").Value = "1.003"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("D42").Value = "0.8121"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "0.0₈125"
$ws.Range("E43").Value = "  +4.32%  "
$ws.Range("D44").Value = "5.440"
$ws.Range("D45").Value = "98.70"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("D46").Value = "1.768.28"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").Value = "0.4521"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").Value = "54.54"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  +0.01%  "

# Restore default styling (drop the temporary text-format override)
# so cells keep their original (unstyled) appearance.
$ws.Range("D2:E51").Style = "Normal"
